$wb = $excel.ActiveWorkbook

# Sheet 1: VENTAS POR GRUPO - MEGAMAFERS S.A. / PORCELANATO (row 12, column M)
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M12").Value = -5.76

# Sheet 2: VENTA MENSUAL - MEGAMAFERS S.A. / julio (row 12, column F) and total (row 22)
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F12").Value = -5.76
$wsMensual.Range("F22").Value = 27931.54

# Sheet 3: CUMPLIMIENTO MENSUAL - PORCELANATO group (row 16) and TOTAL (row 19)
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D16").Value = 19850.12
$wsCumplimiento.Range("E16").Value = 24416.12
$wsCumplimiento.Range("F16").Value = 0.4484257077176648

$wsCumplimiento.Range("D19").Value = 27931.54
$wsCumplimiento.Range("E19").Value = 37446.45762291769
$wsCumplimiento.Range("F19").Value = 0.4272315001309988
